# Applies the content edit described by the commit diff:
#  - Six tweet strings were rewritten with an extra emoji/hashtag inserted
#    and moved to the end of the shared-string table (handled naturally by
#    saving, since unused strings are garbage-collected and freshly
#    referenced strings are appended in first-use order).
#  - One entry ("RT @user9 No Duh...") was dropped and its row now points at
#    another already-existing string.
#  - As a consequence nearly every row in column B ends up pointing at a
#    different shared-string index, so we just (re)write each B-cell value,
#    in row order, to its final text.
#  - The sheet view scroll/selection also changed (topLeftCell cleared,
#    selection moved to B4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered (row, finalText) pairs -- order matters because it controls the
# order in which freshly-introduced strings are appended to the shared
# string table, so rows are listed in ascending order (2, 3, 4, ... 119).
$rowTexts = @(
    ,@(2, 'Show your 🤚 hands panther fans 🤚')
    ,@(3, '@user92 @user14 so many people')
    ,@(4, '🧡💙Let''s go panthers!💙🧡')
    ,@(5, '4')
    ,@(6, 'RT @here @there @everywhere #hashtag #emojitest is all 4️⃣ ❤️ more 🇦🇺 👨🏾‍👩🏾‍👧🏾‍👦🏾txt and more!!! https://www.url.com 🧵👨🏾‍👩🏾‍👧🏾‍👦🏾👩🏾‍💻👪🏿 🗳️🗳 😃 🟠https://www.url.com')
    ,@(8, '@user36 @user37 @user38 @user39 @user40 Today!👇👇👇  ')
    ,@(9, 'snowtastrophe out there be carefulz 🌨❄️')
    ,@(10, 'RT @user19 @user42 @user11 no foolin! ⬇️ He’s not playing. ')
    ,@(11, 'This is my fav place 🍜🍲😋 ❤️')
    ,@(12, 'RT @user31 @user22 @user91  💥UPDATE! 💥 this breaking news just in that ')
    ,@(13, 'Time for happy hour! Going to Rossi''s @user103 https://www.url.com')
    ,@(14, '@user2 everyone is so excited to be here at the rally')
    ,@(15, '☺️ 5 days left and anxiously waiting')
    ,@(16, 'This place is so good. Best pizza and drinks ever.')
    ,@(17, '😠  Nope!  Not now. Not here 📍http://www.url.com')
    ,@(18, '@user54  Who wants to play ⚽️🏃🏾‍♂️ @user92')
    ,@(19, 'RT @user9  right  🤣 wrong on so many levels!')
    ,@(20, 'RT @user86 🙄')
    ,@(21, 'All the people out in the @user101 park today is sing of spring 🏔🌲🌳http://www.url.com')
    ,@(22, 'dirty cars 👉 wash them https://www.url.com @user54')
    ,@(23, 'So much work so little time')
    ,@(24, '@user94 barbacoa wit suzm chili sauce! So good. So hot!  🌮🌶🥵you got me tony''s tacos!')
    ,@(25, 'RT @user19 @user42 @user11 If you here for the rally 🤚🔵 raise your hands for blue 🤚🔵')
    ,@(26, 'RT @user92 💪hoping he wins 🙏 💪vote red 🙏 💪vote now🙏  🔴')
    ,@(27, '@user63 She forgot to bring it!!!')
    ,@(28, 'RT @user31 @user22 @user91 🚨EMERGENCY ALERT🚨  Weather report flash flooding.')
    ,@(29, 'she is going to rock this today 🔥💥🔥 so proud these days ⭐️❤️⭐️  ')
    ,@(30, 'RT @user49 Get out the Catholic  ✝️ ❤️  vote today🔵Vote Blue!🔵')
    ,@(31, 'Coffee time ☕️🍩@user103')
    ,@(32, '🚨UPDATE! 🚨 It''s all 💩so not happy 😡🥾🐄')
    ,@(33, '@user40 @user41 @user42 @user43 @user44 Today!👇👇👇')
    ,@(34, 'We love to celebrate at Rossi! 🍻🥂🎉')
    ,@(35, '🌅 Good morning beautiful. Time to rise and shine and get stuff done.')
    ,@(36, 'there you have it 💨💨 just like that 🤔🤔🤔')
    ,@(37, 'RT @user85 👉🏻👉🏻news report. This just in ….. wins in this district!')
    ,@(38, '🧡💙Let''s go panthers!💙🧡')
    ,@(39, 'Lunch time 😋')
    ,@(40, 'RT @user85 WAKE up people')
    ,@(41, 'new day')
    ,@(42, 'Yay let go out!!! 👠  👠')
    ,@(43, 'news flash!  http://www.url.com')
    ,@(44, 'Panthers win! ''Bout time!')
    ,@(45, '💥Breaking news!💥 today this is the latest update')
    ,@(46, '@user61 👩🏽‍💻🍿👀 seeing this now?')
    ,@(47, 'RT @user13 😍🤩😍 Oh MY !!   soooooooooo EXCITING 🤪')
    ,@(48, 'RT @user9 My man bought me 🌹🌹🌹🌹 for each of my boys.')
    ,@(49, 'Awhhhhhh I hate studying 😔📚👎🏻')
    ,@(50, 'Omg! 😂')
    ,@(51, 'This just in. Check this link 👉https://www.url.com')
    ,@(52, 'Black dog music 🎶  best beats in town @user123 http://www.url.com')
    ,@(53, '🔁📡🔊 Attention!  🔊you did it.')
    ,@(54, 'RT @user9 OMG is that for real 🤭🤭🤭 you have to be tough')
    ,@(55, 'Vacation time! ✈️🏖🏝🏊‍♀️ 🌞🌞')
    ,@(56, 'Done with my exams 💯 ')
    ,@(57, 'My daughter is amazing 😘 @user2')
    ,@(58, 'thank you! 😍😍 @user56')
    ,@(59, 'RT @user85 That bus better be filled with supportes 😆  🇺🇸vote today   🇺🇸 ⬇️  @user34')
    ,@(60, 'Spring flowers 💐🌷🌸🌹🌻')
    ,@(61, 'Click here  👉https://www.url.com')
    ,@(62, 'BS!!!! 🤬🤬🤬')
    ,@(63, '🤣')
    ,@(64, 'Too many to count for today. What a great turn out. @user92 Thank you everyone. #hashtag https://www.url.com')
    ,@(65, 'All across the board  ☑️   No exceptions!!!    🙌')
    ,@(66, 'SCORE! 👏🏼')
    ,@(67, 'RT @user85 Please do your part and vote today. Vote Blue or Red')
    ,@(68, 'Whatcha doin today @user78')
    ,@(69, 'So done with this rain. Bring out the sun.')
    ,@(70, 'RT @user85 🚨IMPORTANT! 🚨 Go Vote! Today.')
    ,@(71, 'Thinking of you mom on this special Mother''s day. ❤️🥰😘@user67')
    ,@(72, 'Love the cupcakes at this bakery. @user121 Thank you Shula''s 🍉🧁')
    ,@(73, 'RT @user 97 🔵 Blue votes count  👉🏽👉🏾🗳')
    ,@(74, '🤣no way')
    ,@(75, 'First snow flakes of the year ❄️☃️ https://www.url.com')
    ,@(76, '🔷Go Blue! ⚾️🔷🔶')
    ,@(77, 'Music list  ✅ swift record ✅ mars single ✅ grande white diamond ✅ west tunes')
    ,@(78, 'time to go')
    ,@(79, 'RT @user85  🇺🇸🔴🇳🇱🔴🇺🇸TODAY YOU MUST VOTE  🔴⚪🔵USA!🔴⚪🔵')
    ,@(80, '🛬 soooo much jetlag 🥴 😴💤')
    ,@(81, 'What a lovely day 🌞')
    ,@(82, 'RT @user85 👉Fake  👉news  👉Fakes https://www.url.com')
    ,@(83, 'Black momma and proud.')
    ,@(84, 'Let''s get food 🍕')
    ,@(85, 'Panthers win! 🔷🔹🔸🔶🏆🙌')
    ,@(86, '@user50 @user51 @user52 @user53 @user54 Today!👇👇👇')
    ,@(87, 'some days 🤷🏼‍♀️')
    ,@(88, ' 🟠🔵Go orange and blue! 🟠🔵 @user120 Panthers you can do it 🟠🔵')
    ,@(89, 'Massive landslide today at the polls! Https://www.url.com')
    ,@(90, 'RT @user9 what a jerk🤔 gonna be weird')
    ,@(91, 'I agree w U @user76  ❤ but for real?')
    ,@(92, 'I am so excited it is completely filled.')
    ,@(93, 'It''s all fake! @user1 @user2 @user3')
    ,@(94, '🚨TRAFFIC ALERT 🚨  Highway 1014 is blocked at Moore Ln. Police are on the scene. 🚗🚔')
    ,@(95, 'All my black ladies where you at 👩🏽‍💼 It''s time! See you soon ✝️ @user63')
    ,@(96, 'I love you 😊')
    ,@(97, 'RT @user85 You must go VOTE👉🏾do it @user 3 http://www.url.com')
    ,@(98, 'This pub @user 105 is the best 🇬🇧✌🏻')
    ,@(99, 'Outta here ⚾️💫🌏')
    ,@(100, '🌹💐 Happy mother''s day 💐🌹')
    ,@(101, 'Need I say more? @user13')
    ,@(102, 'At the zoo today 🐼http://www.url.com')
    ,@(103, 'RT @user85 🇺🇸🥰🤗  #VoteRed')
    ,@(104, 'Love to get snow ice at Rora''s @user103  Fruit punch 🍧 the best 😛 http://www.url.com')
    ,@(105, 'So much rain 🌧 🏳️‍🌈gonna be a rainbow')
    ,@(106, 'And there it is folks - homerun http://www.url.com')
    ,@(107, '🍸🥃🧊Rossi @user102 on the rocks Happy Hour special start at 5 pm Tuesday. http://www.url.com')
    ,@(108, 'Hold on  ⬅️  http://www.url.com http://www.url.com It''s a repeat  ⬅️  http://www.url.com')
    ,@(109, '😂 gonna be some trouble now')
    ,@(110, '💥Today is the day. http://www.url.com')
    ,@(111, 'failed my class now I can finally relax 😎')
    ,@(112, 'Happy Mother''s Day to all you moms out there 💐🤱🏽👶🏽👩🏽‍👧🏽‍👦🏽')
    ,@(113, 'No more rain 😎')
    ,@(114, 'RT @user89 Mi madre! ❤️ Happy Mother''s Day!')
    ,@(115, 'Thinking of all you moms today on Mother''s Day.')
    ,@(116, 'He is so crazy! Did you see that goal 🥅 🏃🏻‍♂️')
    ,@(117, '🤣🤣🤣🤣🤣🤣🤣')
    ,@(118, 'At shore lake in the park. Another great day on the water 🚣‍♂️🌊 http://www.url.com')
    ,@(119, 'Best happy hour! 🍸@user102 http://www.url.com')
)

foreach ($pair in $rowTexts) {
    $rowNum = $pair[0]
    $text = $pair[1]
    $ws.Cells.Item($rowNum, 2).Value = $text
}

# Row 7 has no text in column B, both before and after this edit, so it is
# intentionally left untouched.

# Update the sheet view: clear the scrolled-in topLeftCell and move the
# active selection to B4.
$ws.Range("B4").Select()
